# feat: add 2022-Q1 data
#
# - Inserts a new worksheet "2022-Q1" right after "2021-Q4" (and right
#   before "总计"), populated with the per-fund holdings for that quarter.
# - Appends a matching summary row to the "总计" sheet (inserted as the
#   new first data row, since it is sorted most-recent-first).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Helpers
# ---------------------------------------------------------------------------

function Set-HeaderCell($range) {
    # Mirrors the workbook's existing "header / row-index" cell style:
    # bold, thin box border, centered horizontally, top-aligned vertically.
    $range.Font.Bold = $true
    $range.HorizontalAlignment = -4108   # xlCenter
    $range.VerticalAlignment = -4160     # xlTop
    $range.Borders.LineStyle = 1         # xlContinuous
    $range.Borders.Weight = 2            # xlThin
}

function Set-TextValue($range, [string]$text) {
    # Force the cell to be stored as text even though it looks numeric
    # (matches the source data, e.g. "84.75" kept verbatim as a string).
    $range.NumberFormat = "@"
    $range.Value = $text
}

# ---------------------------------------------------------------------------
# 1) Insert the new "2022-Q1" sheet right after "2021-Q4"
# ---------------------------------------------------------------------------

$q4sheet = $wb.Worksheets.Item("2021-Q4")
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $q4sheet)
$newSheet.Name = "2022-Q1"

# Header row
$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $col = $i + 2   # headers start at column B
    $cell = $newSheet.Cells.Item(1, $col)
    $cell.Value = $headers[$i]
    Set-HeaderCell $cell
}

# Fund rows: A, code, name, scale(text), position(text), ratio(text), value(text), rank(number)
$fundRows = @(
    @(0, "008545", "泓德丰润三年持有期混合", "84.75", "91.10", "5.19", "4.3985", 7),
    @(1, "005395", "泓德臻远回报灵活配置混合", "33.94", "93.62", "6.08", "2.0636", 5),
    @(2, "001500", "泓德远见回报混合",        "26.71", "93.72", "6.27", "1.6747", 9),
    @(3, "004965", "泓德致远混合A",           "21.03", "46.32", "5.31", "1.1167", 3),
    @(4, "004966", "泓德致远混合C",           "2.92",  "46.32", "5.31", "0.1551", 3),
    @(5, "255010", "国联安稳健混合",          "2.33",  "69.07", "3.49", "0.0813", 10),
    @(6, "006863", "国联安智能制造混合",      "0.23",  "93.60", "5.52", "0.0127", 9)
)

$rowNum = 2
foreach ($fr in $fundRows) {
    $aCell = $newSheet.Cells.Item($rowNum, 1)
    $aCell.Value = $fr[0]
    Set-HeaderCell $aCell

    # Fund code (e.g. "008545") must stay text so leading zeros survive.
    Set-TextValue ($newSheet.Cells.Item($rowNum, 2)) $fr[1]

    $cCell = $newSheet.Cells.Item($rowNum, 3)
    $cCell.Value = $fr[2]

    Set-TextValue ($newSheet.Cells.Item($rowNum, 4)) $fr[3]
    Set-TextValue ($newSheet.Cells.Item($rowNum, 5)) $fr[4]
    Set-TextValue ($newSheet.Cells.Item($rowNum, 6)) $fr[5]
    Set-TextValue ($newSheet.Cells.Item($rowNum, 7)) $fr[6]

    $hCell = $newSheet.Cells.Item($rowNum, 8)
    $hCell.Value = $fr[7]

    $rowNum++
}

$newSheet.Range("A1:H8").Columns.AutoFit() | Out-Null

# ---------------------------------------------------------------------------
# 2) Prepend the new quarter to the "总计" (totals) sheet
# ---------------------------------------------------------------------------

$totalSheet = $wb.Worksheets.Item("总计")

# Existing data lives in rows 2-6 (2021-Q4 .. 2020-Q4). Push it all down one
# row to make room for the new 2022-Q1 row, then rewrite column A's running
# index (0,1,2,...) across the whole block so it stays contiguous.
$totalSheet.Rows.Item(2).Insert()

$totalRows = @(
    @(0, "2022-Q1", 7,  9.5),
    @(1, "2021-Q4", 9,  16.62),
    @(2, "2021-Q3", 9,  18.26),
    @(3, "2021-Q2", 8,  17.18),
    @(4, "2021-Q1", 12, 23.9),
    @(5, "2020-Q4", 9,  19.11)
)

$rowNum = 2
foreach ($tr in $totalRows) {
    $aCell = $totalSheet.Cells.Item($rowNum, 1)
    $aCell.Value = $tr[0]
    Set-HeaderCell $aCell

    $totalSheet.Cells.Item($rowNum, 2).Value = $tr[1]
    $totalSheet.Cells.Item($rowNum, 3).Value = $tr[2]
    $totalSheet.Cells.Item($rowNum, 4).Value = $tr[3]

    $rowNum++
}

$totalSheet.Range("A1:D7").Columns.AutoFit() | Out-Null
